# Applies the "added likes functionality + sorting" edit described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31: add hours worked (D31) and task description (E31) for date 2022-05-16 (44681)
$ws.Range("D31").Value = 5
$ws.Range("E31").Value = "forbedret søk og sortering. Implementert liking av innhold"

# Widen column E to fit the new (longer) text, matching the bestFit width recorded in the file
$ws.Columns.Item(5).ColumnWidth = 53

# Update the selected cell to reflect where the author ended up after editing
$ws.Range("H32").Select()

$wb.Save()
